# Applies the "Updated Data and Main" commit to Config.xlsx
#  - Settings!B2/C2: queue name + description updated (ProcessABCQueue -> EmailsQueue)
#  - Constants: five new Email* rows (12-16) + a stray hyperlink-styled blank cell (C17)
#  - selections / active sheet updated to match the new editing session
#  - row heights tweak (30 -> 28.5) to match the re-saved layout

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("B2").Value = "EmailsQueue"
$settings.Range("C2").Value = "Orchestrator queue Name. Ensure the queue unique property is not checked."

$settings.Rows.Item(4).RowHeight = 28.5

$settings.Activate()
$settings.Range("C2").Select()

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

$constants.Rows.Item(2).RowHeight = 28.5

$constants.Range("A12").Value = "EmailsFilePath"
$constants.Range("B12").Value = "Data\Input\Emails.xlsx"
$constants.Range("C12").Value = "Relative file path to the Emails.xlsx file"

$constants.Range("A13").Value = "EmailAccount"
$constants.Range("B13").Value = "MPLURAD1@depaul.edu"
$constants.Range("C13").Value = 'Email address for "Account" property within send outlook mail message activity'

$constants.Range("A14").Value = "EmailSenderName"
$constants.Range("B14").Value = "Miles Plurad"
$constants.Range("C14").Value = 'Sender name for "SendOnBehalfOfName" property within send outlook mail message activity'

$constants.Range("A15").Value = "EmailSubject"
$constants.Range("B15").Value = "P2 and P3 Invitations"
$constants.Range("C15").Value = 'Email subject for "Subject" property within send outlook mail message activity'

$constants.Range("A16").Value = "EmailBody"
$constants.Range("B16").Value = "Here is the link to P2 and P3 presentations: https://github.com/mplurad/uipath-automation-10"
$constants.Range("C16").Value = 'Email message for "Body" property within send outlook mail message activity'

# C17 picked up the built-in "Hyperlink" style (underlined theme-10 font) from a
# link that was added and then cleared, leaving just the formatting behind.
$constants.Hyperlinks.Add($constants.Range("C17"), "https://github.com/mplurad/uipath-automation-10")
$constants.Range("C17").Hyperlinks.Delete()
$constants.Range("C17").ClearContents()

$constants.Activate()
$constants.Range("B16").Select()

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Activate()
$assets.Range("D7").Select()

# ---------------------------------------------------------------------------
# Leave "Constants" as the active tab, matching the saved workbook view.
# ---------------------------------------------------------------------------
$constants.Activate()
